# Update the "取得日時" (acquired timestamp) column on the active sheet
# ("ランサーズ") from the previous run's timestamp to the latest run's
# timestamp, reflecting a re-run of the scraper at 2025-09-12 18:28:31.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-12 18:28:31"

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $newTimestamp
}
